$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections ---
$ws.Range("E5").Value = "fecha_registro"
$ws.Range("A1").Value = "Proveedores"

# --- Highlight (yellow fill) the header blocks for Proveedores, Productos and Usuarios tables ---
$yellow = 65535  # RGB(255,255,0)

foreach ($addr in @("A1:F1","A2:F2","A4:F4","A5:F5","A10:F10","A11:F11")) {
    $ws.Range($addr).Interior.Color = $yellow
}

# --- Sheet view changes ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 113
$ws.Range("E12").Select()
